# Lab01 Review Report - translate header/labels to German and widen column B
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Titel: BeCruel (Mobile Project)"
$ws.Range("A6").Value = "Datum: 10.01.2023"
$ws.Range("A7").Value = "Verbesserungen: "

$ws.Range("B8").Value = "Klassennamen umformen für bessere Liesbarkeit des Kodes"
$ws.Range("B9").Value = "Mehr Fehlerhandlung für Corner Cases"
$ws.Range("B10").Value = "Branch-löschen um Redundanz zu eliminieren um Nutzbarkeit zu verbessern/Kode-Lisibilität-Niveau zu erhöhen"
$ws.Range("B11").Value = "Verändern des Output Formatierung  für einige der Features"

# Widen column B to fit the longer German text
$ws.Columns.Item(2).ColumnWidth = 49.25

# Move the active selection to B12, as last left by the editing user
$ws.Range("B12").Select()
